$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), row 15
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 2158.8462
$ws.Cells.Item(15, 9).Value = 2158.8462
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 6476.5386
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -6307.5386

# Sheet ALC (index 1), row 33
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(33, 8).Value = 367
$ws.Cells.Item(33, 9).Value = 367
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 367
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -138

# Sheet ALC (index 1), row 70
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(70, 8).Value = 3355.0557
$ws.Cells.Item(70, 9).Value = 2279.4
$ws.Cells.Item(70, 10).Value = 3768.7693
$ws.Cells.Item(70, 11).Value = 6838.200000000001
$ws.Cells.Item(70, 12).Value = 11306.3079
$ws.Cells.Item(70, 13).Value = -6568.200000000001
$ws.Cells.Item(70, 14).Value = -11846.3079

# Sheet ALC (index 1), row 73
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(73, 8).Value = 3355.0557
$ws.Cells.Item(73, 9).Value = 2279.4
$ws.Cells.Item(73, 10).Value = 3768.7693
$ws.Cells.Item(73, 11).Value = 6838.200000000001
$ws.Cells.Item(73, 12).Value = 11306.3079
$ws.Cells.Item(73, 13).Value = -5902.200000000001
$ws.Cells.Item(73, 14).Value = -13178.3079

# Sheet ARM (index 2), row 32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 4785.8623
$ws.Cells.Item(32, 9).Value = 4785.8623
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 4785.8623
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -4498.8623

# Sheet ARM (index 2), row 35
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(35, 8).Value = 2158.2
$ws.Cells.Item(35, 9).Value = 833.3333
$ws.Cells.Item(35, 10).Value = 4145.5
$ws.Cells.Item(35, 11).Value = 833.3333
$ws.Cells.Item(35, 12).Value = 4145.5
$ws.Cells.Item(35, 13).Value = -427.3333
$ws.Cells.Item(35, 14).Value = -4957.5

# Sheet ARM (index 2), row 61
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 4332
$ws.Cells.Item(61, 9).Value = 3999
$ws.Cells.Item(61, 10).Value = 4498.5
$ws.Cells.Item(61, 11).Value = 3999
$ws.Cells.Item(61, 12).Value = 4498.5
$ws.Cells.Item(61, 13).Value = -3787
$ws.Cells.Item(61, 14).Value = -4922.5

# Sheet ARM (index 2), row 132
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(132, 8).Value = 1900
$ws.Cells.Item(132, 9).Value = 1125
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 3375
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -845
$ws.Cells.Item(132, 14).Value = -20060

# Sheet ARM (index 2), row 136
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(136, 8).Value = 4332
$ws.Cells.Item(136, 9).Value = 3999
$ws.Cells.Item(136, 10).Value = 4498.5
$ws.Cells.Item(136, 11).Value = 11997
$ws.Cells.Item(136, 12).Value = 13495.5
$ws.Cells.Item(136, 13).Value = -9447
$ws.Cells.Item(136, 14).Value = -18595.5

# Sheet BSM (index 3), row 20
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 714
$ws.Cells.Item(20, 9).Value = 714
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 714
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -467
$ws.Cells.Item(20, 14).Value = ""

# Sheet BSM (index 3), row 105
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(105, 8).Value = 1222.1428
$ws.Cells.Item(105, 9).Value = 1000
$ws.Cells.Item(105, 10).Value = 1259.1666
$ws.Cells.Item(105, 11).Value = 1000
$ws.Cells.Item(105, 12).Value = 1259.1666
$ws.Cells.Item(105, 13).Value = 747
$ws.Cells.Item(105, 14).Value = -4753.1666

# Sheet BSM (index 3), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 4025.0881
$ws.Cells.Item(134, 9).Value = 4164
$ws.Cells.Item(134, 10).Value = 3219.4
$ws.Cells.Item(134, 11).Value = 12492
$ws.Cells.Item(134, 12).Value = 9658.200000000001
$ws.Cells.Item(134, 13).Value = -9957
$ws.Cells.Item(134, 14).Value = -14728.2

# Sheet BSM (index 3), row 137
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(137, 8).Value = 46999.2
$ws.Cells.Item(137, 9).Value = 35000
$ws.Cells.Item(137, 10).Value = 49999
$ws.Cells.Item(137, 11).Value = 35000
$ws.Cells.Item(137, 12).Value = 49999
$ws.Cells.Item(137, 13).Value = -29900
$ws.Cells.Item(137, 14).Value = -60199

# Sheet CRP (index 4), row 6
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 8).Value = 167191.67
$ws.Cells.Item(6, 9).Value = 250237.75
$ws.Cells.Item(6, 10).Value = 1099.5
$ws.Cells.Item(6, 11).Value = 250237.75
$ws.Cells.Item(6, 12).Value = 1099.5
$ws.Cells.Item(6, 13).Value = -250124.75
$ws.Cells.Item(6, 14).Value = -1325.5

# Sheet CRP (index 4), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3268.88
$ws.Cells.Item(31, 9).Value = 2301.0476
$ws.Cells.Item(31, 10).Value = 8350
$ws.Cells.Item(31, 11).Value = 2301.0476
$ws.Cells.Item(31, 12).Value = 8350
$ws.Cells.Item(31, 13).Value = -2006.0476
$ws.Cells.Item(31, 14).Value = -8940

# Sheet CRP (index 4), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 3268.88
$ws.Cells.Item(34, 9).Value = 2301.0476
$ws.Cells.Item(34, 10).Value = 8350
$ws.Cells.Item(34, 11).Value = 2301.0476
$ws.Cells.Item(34, 12).Value = 8350
$ws.Cells.Item(34, 13).Value = -2099.0476
$ws.Cells.Item(34, 14).Value = -8754

# Sheet CRP (index 4), row 99
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(99, 8).Value = 4799.7144
$ws.Cells.Item(99, 9).Value = 5499.5
$ws.Cells.Item(99, 10).Value = 3866.6667
$ws.Cells.Item(99, 11).Value = 5499.5
$ws.Cells.Item(99, 12).Value = 3866.6667
$ws.Cells.Item(99, 13).Value = -4001.5
$ws.Cells.Item(99, 14).Value = -6862.6667

# Sheet CRP (index 4), row 122
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(122, 8).Value = 2502.2
$ws.Cells.Item(122, 9).Value = 2252.75
$ws.Cells.Item(122, 10).Value = 3500
$ws.Cells.Item(122, 11).Value = 6758.25
$ws.Cells.Item(122, 12).Value = 10500
$ws.Cells.Item(122, 13).Value = -4308.25
$ws.Cells.Item(122, 14).Value = -15400

# Sheet CRP (index 4), row 126
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(126, 8).Value = 4799.7144
$ws.Cells.Item(126, 9).Value = 5499.5
$ws.Cells.Item(126, 10).Value = 3866.6667
$ws.Cells.Item(126, 11).Value = 16498.5
$ws.Cells.Item(126, 12).Value = 11600.0001
$ws.Cells.Item(126, 13).Value = -14028.5
$ws.Cells.Item(126, 14).Value = -16540.0001

# Sheet CRP (index 4), row 134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 2039.25
$ws.Cells.Item(134, 9).Value = 1909.9131
$ws.Cells.Item(134, 10).Value = 5014
$ws.Cells.Item(134, 11).Value = 5729.7393
$ws.Cells.Item(134, 12).Value = 15042
$ws.Cells.Item(134, 13).Value = -3194.7393
$ws.Cells.Item(134, 14).Value = -20112

# Sheet CUL (index 5), row 75
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(75, 8).Value = 4325.75
$ws.Cells.Item(75, 9).Value = 3832
$ws.Cells.Item(75, 10).Value = 4622
$ws.Cells.Item(75, 11).Value = 11496
$ws.Cells.Item(75, 12).Value = 13866
$ws.Cells.Item(75, 13).Value = -10498
$ws.Cells.Item(75, 14).Value = -15862

# Sheet CUL (index 5), row 78
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(78, 8).Value = 4325.75
$ws.Cells.Item(78, 9).Value = 3832
$ws.Cells.Item(78, 10).Value = 4622
$ws.Cells.Item(78, 11).Value = 34488
$ws.Cells.Item(78, 12).Value = 41598
$ws.Cells.Item(78, 13).Value = -29496
$ws.Cells.Item(78, 14).Value = -51582

# Sheet CUL (index 5), row 132
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(132, 8).Value = 1378.8
$ws.Cells.Item(132, 9).Value = 447
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 4023
$ws.Cells.Item(132, 12).Value = 18000
$ws.Cells.Item(132, 13).Value = -1493
$ws.Cells.Item(132, 14).Value = -23060

# Sheet CUL (index 5), row 140
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(140, 8).Value = 1768.1666
$ws.Cells.Item(140, 9).Value = 1428.909
$ws.Cells.Item(140, 10).Value = 5500
$ws.Cells.Item(140, 11).Value = 4286.727000000001
$ws.Cells.Item(140, 12).Value = 16500
$ws.Cells.Item(140, 13).Value = 893.2729999999992
$ws.Cells.Item(140, 14).Value = -26860

# Sheet GSM (index 6), row 2
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 300.15
$ws.Cells.Item(2, 9).Value = 345.06668
$ws.Cells.Item(2, 10).Value = 165.4
$ws.Cells.Item(2, 11).Value = 345.06668
$ws.Cells.Item(2, 12).Value = 165.4
$ws.Cells.Item(2, 13).Value = -232.06668
$ws.Cells.Item(2, 14).Value = -391.4

# Sheet GSM (index 6), row 19
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(19, 8).Value = 509.7143
$ws.Cells.Item(19, 9).Value = 594.5
$ws.Cells.Item(19, 10).Value = 1
$ws.Cells.Item(19, 11).Value = 594.5
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = -306.5
$ws.Cells.Item(19, 14).Value = -577

# Sheet GSM (index 6), row 132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 2223.6428
$ws.Cells.Item(132, 9).Value = 2368.1667
$ws.Cells.Item(132, 10).Value = 1356.5
$ws.Cells.Item(132, 11).Value = 7104.500100000001
$ws.Cells.Item(132, 12).Value = 4069.5
$ws.Cells.Item(132, 13).Value = -4574.500100000001
$ws.Cells.Item(132, 14).Value = -9129.5

# Sheet LTW (index 7), row 7
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4490.2
$ws.Cells.Item(7, 9).Value = 4362.75
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 4362.75
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -4250.75
$ws.Cells.Item(7, 14).Value = -5224

# Sheet LTW (index 7), row 126
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(126, 8).Value = 4490.2
$ws.Cells.Item(126, 9).Value = 4362.75
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 13088.25
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = -10618.25
$ws.Cells.Item(126, 14).Value = -19940

# Sheet WVR (index 8), row 54
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(54, 8).Value = 24714.285
$ws.Cells.Item(54, 9).Value = 12000
$ws.Cells.Item(54, 10).Value = 29800
$ws.Cells.Item(54, 11).Value = 12000
$ws.Cells.Item(54, 12).Value = 29800
$ws.Cells.Item(54, 13).Value = -11480
$ws.Cells.Item(54, 14).Value = -30840
